# backdatedVinTable_UT_SS.xlsx edit
#
# The sheet had two columns that are no longer needed: column Z ("STAT")
# and what was originally column AC ("CHOICE_TIER"). Removing them shifts
# every later column left (AA3 comment -> Z3, AL -> AJ, etc.) and drops the
# now-unused shared strings (STAT / CHOICE_TIER / S / RT) automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The cell comment is anchored to AA3 ("COLL_SYMBOL" header column); once
# column Z is removed that same comment needs to live on Z3. Capture its
# text first so it can be re-created at the new location.
$cmt = $ws.Range("AA3").Comment
$cmtText = $cmt.Text()
$cmt.Delete()

# Delete column Z ("STAT"). Everything to its right shifts one column left,
# so the old "CHOICE_TIER" column (AC) is now AB.
$ws.Columns("Z").Delete()
$ws.Columns("AB").Delete()

$newCmt = $ws.Range("Z3").AddComment($cmtText)

# Restore the view state recorded in the saved workbook: scrolled so column
# T is the left-most visible column, with AA8 selected (below/right of the
# populated data, matching the post-edit selection).
$ws.Range("AA8").Select()
$excel.ActiveWindow.ScrollColumn = 20
